$d = $word.ActiveDocument

# Step 1: remove "snd ^ " prefix (before runAssumingCascade)
$d.Content.Find.Execute("snd ^ runAssumingCascade", $true, $false, $false, $false, $false, $true, 1, $false, "runAssumingCascade", 2) | Out-Null

# Step 2: append " |> snd" after "chain observation world"
$d.Content.Find.Execute("chain observation world", $true, $false, $false, $false, $false, $true, 1, $false, "chain observation world |> snd", 2) | Out-Null
